# Updated RAD Test Scripts and Test Data for Existing Liability.
#
# "Existing Liability w/Notice Number" is renamed to the fuller
# "Existing Liability with Notice/Invoice Number" (rows 2 and 6), and
# every data row now has its Execute flag ("Y") filled in (rows 2, 3, 5,
# 6, 7 previously had a blank Execute cell; row 4 already had "Y").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPaymentType = "Existing Liability with Notice/Invoice Number"

# Row 2: Existing Liability w/Notice Number, Personal Income Tax
$ws.Range("C2").Value = "Y"
$ws.Range("D2").Value = $newPaymentType

# Row 3: Quarterly Estimated Tax, Personal Income Tax
$ws.Range("C3").Value = "Y"

# Row 4: Extension Payments, Personal Income Tax (Execute already "Y")
$ws.Range("C4").Value = "Y"

# Row 5: New Tax Return Amount Due, Personal Income Tax
$ws.Range("C5").Value = "Y"

# Row 6: Existing Liability w/Notice Number, Estate Tax
$ws.Range("C6").Value = "Y"
$ws.Range("D6").Value = $newPaymentType

# Row 7: New Tax Return Amount Due, Estate Tax
$ws.Range("C7").Value = "Y"

# Match the selection left behind in the saved workbook (Execute column
# C2:C7, active cell C2).
$ws.Range("C2:C7").Select()

# The Date column (bestFit) was re-measured by Excel on save; nudge its
# stored width as close as this engine's pixel-snapped ColumnWidth allows.
$ws.Columns.Item(2).ColumnWidth = 26.8
